# Update UML diagrams and corresponding images
#
# 1) The "date last saved" auto-date field shown on the slide master and
#    every slide layout moved on from 5/3/2018 to 8/3/2018.
# 2) The variable name shown in the undo/redo stack table was renamed
#    from prevTaskBook to prevOrganizer.

$p = $ppt.ActivePresentation

# --- 1. Refresh the "datetimeFigureOut" placeholder text wherever it
#        shows the old date, across the slide master and all of its
#        custom (slide) layouts. ---
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
  $sh = $m.Shapes.Item($i)
  if ($sh.HasTextFrame) {
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -eq "5/3/2018") {
      $tr.Text = "8/3/2018"
    }
  }
}

$layouts = $m.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
  $cl = $layouts.Item($j)
  for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
    $sh = $cl.Shapes.Item($i)
    if ($sh.HasTextFrame) {
      $tr = $sh.TextFrame.TextRange
      if ($tr.Text -eq "5/3/2018") {
        $tr.Text = "8/3/2018"
      }
    }
  }
}

# --- 2. Rename prevTaskBook -> prevOrganizer inside the "Table 20"
#        table on slide 1 (second row, first/only column, second
#        paragraph). ---
$s = $p.Slides.Item(1)
$tableShape = $s.Shapes.Item("Table 20")
$tbl = $tableShape.Table
$cell = $tbl.Cell(2, 1)
$para = $cell.Shape.TextFrame.TextRange.Paragraphs(2)
$para.Text = $para.Text.Replace("prevTaskBook", "prevOrganizer")
